$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 59475.5
$ws.Range("J17").Value = 59475.5
$ws.Range("L17").Value = 178426.5
$ws.Range("N17").Value = -178762.5
$ws.Range("H19").Value = 2977.7778
$ws.Range("I19").Value = 209.16667
$ws.Range("J19").Value = 4362.0835
$ws.Range("K19").Value = 209.16667
$ws.Range("L19").Value = 4362.0835
$ws.Range("M19").Value = -34.16667000000001
$ws.Range("N19").Value = -4712.0835
$ws.Range("H29").Value = 3025.75
$ws.Range("I29").Value = 3025.75
$ws.Range("K29").Value = 9077.25
$ws.Range("M29").Value = -8796.25
$ws.Range("H46").Value = 1796.9
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1796.9
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 5390.700000000001
$ws.Range("N46").Value = -5628.700000000001
$ws.Range("H60").Value = 1796.9
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1796.9
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").Value = 5390.700000000001
$ws.Range("N60").Value = -6358.700000000001
$ws.Range("H70").Value = 3162.9167
$ws.Range("I70").Value = 2317
$ws.Range("J70").Value = 3283.762
$ws.Range("K70").Value = 6951
$ws.Range("L70").Value = 9851.286
$ws.Range("M70").Value = -6681
$ws.Range("N70").Value = -10391.286
$ws.Range("H73").Value = 3162.9167
$ws.Range("I73").Value = 2317
$ws.Range("J73").Value = 3283.762
$ws.Range("K73").Value = 6951
$ws.Range("L73").Value = 9851.286
$ws.Range("M73").Value = -6015
$ws.Range("N73").Value = -11723.286

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 535.4103
$ws.Range("I74").Value = 483.3
$ws.Range("J74").Value = 709.1111
$ws.Range("K74").Value = 483.3
$ws.Range("L74").Value = 709.1111
$ws.Range("M74").Value = 390.7
$ws.Range("N74").Value = -2457.1111
$ws.Range("H77").Value = 535.4103
$ws.Range("I77").Value = 483.3
$ws.Range("J77").Value = 709.1111
$ws.Range("K77").Value = 2416.5
$ws.Range("L77").Value = 3545.5555
$ws.Range("M77").Value = 1951.5
$ws.Range("N77").Value = -12281.5555

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3712
$ws.Range("I134").Value = 1859.2
$ws.Range("J134").Value = 6800
$ws.Range("K134").Value = 5577.6
$ws.Range("L134").Value = 20400
$ws.Range("M134").Value = -3042.6
$ws.Range("N134").Value = -25470

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1856.0785
$ws.Range("I31").Value = 1092.4572
$ws.Range("J31").Value = 3526.5
$ws.Range("K31").Value = 1092.4572
$ws.Range("L31").Value = 3526.5
$ws.Range("M31").Value = -797.4572000000001
$ws.Range("N31").Value = -4116.5
$ws.Range("H34").Value = 1856.0785
$ws.Range("I34").Value = 1092.4572
$ws.Range("J34").Value = 3526.5
$ws.Range("K34").Value = 1092.4572
$ws.Range("L34").Value = 3526.5
$ws.Range("M34").Value = -890.4572000000001
$ws.Range("N34").Value = -3930.5
$ws.Range("H99").Value = 1998.4615
$ws.Range("I99").Value = 1484.4445
$ws.Range("J99").Value = 2270.5881
$ws.Range("K99").Value = 1484.4445
$ws.Range("L99").Value = 2270.5881
$ws.Range("M99").Value = 13.55549999999994
$ws.Range("N99").Value = -5266.5881
$ws.Range("H105").Value = 2781.3809
$ws.Range("I105").Value = 2495.2104
$ws.Range("K105").Value = 2495.2104
$ws.Range("M105").Value = -748.2103999999999
$ws.Range("H126").Value = 1998.4615
$ws.Range("I126").Value = 1484.4445
$ws.Range("J126").Value = 2270.5881
$ws.Range("K126").Value = 4453.333500000001
$ws.Range("L126").Value = 6811.7643
$ws.Range("M126").Value = -1983.333500000001
$ws.Range("N126").Value = -11751.7643
$ws.Range("H132").Value = 3049.0857
$ws.Range("I132").Value = 2046.1
$ws.Range("J132").Value = 4386.4
$ws.Range("K132").Value = 6138.299999999999
$ws.Range("L132").Value = 13159.2
$ws.Range("M132").Value = -3608.299999999999
$ws.Range("N132").Value = -18219.2
$ws.Range("H134").Value = 2157.7666
$ws.Range("I134").Value = 1409.9048
$ws.Range("J134").Value = 3902.7778
$ws.Range("K134").Value = 4229.7144
$ws.Range("L134").Value = 11708.3334
$ws.Range("M134").Value = -1694.7144
$ws.Range("N134").Value = -16778.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2214.6047
$ws.Range("I68").Value = 724.375
$ws.Range("J68").Value = 2555.2285
$ws.Range("K68").Value = 2173.125
$ws.Range("L68").Value = 7665.685500000001
$ws.Range("M68").Value = -1362.125
$ws.Range("N68").Value = -9287.6855
$ws.Range("H71").Value = 2214.6047
$ws.Range("I71").Value = 724.375
$ws.Range("J71").Value = 2555.2285
$ws.Range("K71").Value = 6519.375
$ws.Range("L71").Value = 22997.0565
$ws.Range("M71").Value = -2463.375
$ws.Range("N71").Value = -31109.0565
$ws.Range("H131").Value = 1893.25
$ws.Range("I131").Value = 1952.375
$ws.Range("J131").Value = 1775
$ws.Range("K131").Value = 5857.125
$ws.Range("L131").Value = 5325
$ws.Range("M131").Value = -817.125
$ws.Range("N131").Value = -15405
$ws.Range("H132").Value = 2572.8
$ws.Range("I132").Value = 1224.625
$ws.Range("J132").Value = 4113.5713
$ws.Range("K132").Value = 11021.625
$ws.Range("L132").Value = 37022.14169999999
$ws.Range("M132").Value = -8491.625
$ws.Range("N132").Value = -42082.14169999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 23650.188
$ws.Range("J134").Value = 21873.8
$ws.Range("L134").Value = 65621.39999999999
$ws.Range("N134").Value = -70691.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2309.2856
$ws.Range("I7").Value = 1673.5
$ws.Range("J7").Value = 2563.6
$ws.Range("K7").Value = 1673.5
$ws.Range("L7").Value = 2563.6
$ws.Range("M7").Value = -1561.5
$ws.Range("N7").Value = -2787.6
$ws.Range("H126").Value = 2309.2856
$ws.Range("I126").Value = 1673.5
$ws.Range("J126").Value = 2563.6
$ws.Range("K126").Value = 5020.5
$ws.Range("L126").Value = 7690.799999999999
$ws.Range("M126").Value = -2550.5
$ws.Range("N126").Value = -12630.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1166.6666
$ws.Range("I81").Value = 833.3333
$ws.Range("J81").Value = 1333.3334
$ws.Range("K81").Value = 1666.6666
$ws.Range("L81").Value = 2666.6668
$ws.Range("M81").Value = -605.6666
$ws.Range("N81").Value = -4788.6668
$ws.Range("H84").Value = 1166.6666
$ws.Range("I84").Value = 833.3333
$ws.Range("J84").Value = 1333.3334
$ws.Range("K84").Value = 8333.333000000001
$ws.Range("L84").Value = 13333.334
$ws.Range("M84").Value = -3029.333000000001
$ws.Range("N84").Value = -23941.334
